$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G3").Value = 79913400000000
$ws.Range("G6").Value = 299383815950000
$ws.Range("G7").Value = 34893020155614.79
$ws.Range("G8").Value = 35509789088895.52
